# Rename products/industries in model aggregation
#
# The underlying shared-string table used to hold generic placeholder
# codes (i001, i002, ... / p001, p002, ...) for the elasTRADE / elasPROD /
# FPROD lookup sheets. This commit swaps those placeholders for the real
# product (pXXXX) and industry (iXXXX) mnemonics, and leaves the active
# sheet on "elasFU" instead of "FPROD".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# elasTRADE (sheet2): column A, rows 2-36 -> product codes
# ---------------------------------------------------------------------
$wsTrade = $wb.Worksheets.Item("elasTRADE")
$productCodes = @(
    "pPARI", "pWHEA", "pOCER", "pFVEG", "pOILS",
    "pSUGB", "pFIBR", "pOTHC", "pANIM", "pFORE",
    "pFISH", "pFOSM", "pOTHM", "pFBTO", "pTXWO",
    "pCOKE", "pREFN", "pCHEM", "pRUBP", "pNMMP",
    "pMETP", "pELEC", "pMACH", "pELCF", "pELCG",
    "pTRDI", "pHWAT", "pWATR", "pCONS", "pTRAD",
    "pHORE", "pTRAN", "pREBA", "pPUBO", "pWAST"
)
for ($i = 0; $i -lt $productCodes.Length; $i++) {
    # Leading apostrophe forces text entry (keeps the existing
    # quote-prefix cell style instead of resetting it to General).
    $wsTrade.Cells.Item($i + 2, 1).Value = "'" + $productCodes[$i]
}

# ---------------------------------------------------------------------
# elasPROD (sheet3) and FPROD (sheet4): column A, rows 2-36 -> industry codes
# ---------------------------------------------------------------------
$industryCodes = @(
    "iPARI", "iWHEA", "iOCER", "iFVEG", "iOILS",
    "iSUGB", "iFIBR", "iOTHC", "iANIM", "iFORE",
    "iFISH", "iFOSM", "iOTHM", "iFBTO", "iTXWO",
    "iCOKE", "iREFN", "iCHEM", "iRUBP", "iNMMP",
    "iMETP", "iELEC", "iMACH", "iELCF", "iELCG",
    "iTRDI", "iHWAT", "iWATR", "iCONS", "iTRAD",
    "iHORE", "iTRAN", "iREBA", "iPUBO", "iWAST"
)

$wsProd = $wb.Worksheets.Item("elasPROD")
for ($i = 0; $i -lt $industryCodes.Length; $i++) {
    $wsProd.Cells.Item($i + 2, 1).Value = "'" + $industryCodes[$i]
}

$wsFprod = $wb.Worksheets.Item("FPROD")
for ($i = 0; $i -lt $industryCodes.Length; $i++) {
    $wsFprod.Cells.Item($i + 2, 1).Value = "'" + $industryCodes[$i]
}

# ---------------------------------------------------------------------
# Move the active tab from FPROD to elasFU
# ---------------------------------------------------------------------
$wsFu = $wb.Worksheets.Item("elasFU")
$wsFu.Select()
